# Apply the Nov 17 2023 GitHub Actions crypto-list refresh: updated
# Price (D) / Volume(1h) (E) figures for most rows, plus a rank swap
# between Cronos and TrustWalletToken (rows 41-42, including their
# Coin name / Link columns).
#
# D/E columns hold free-form text (e.g. "36.364.88", "  -2.64%  ") that
# must stay text, not get reinterpreted as numbers by Excel. Force the
# cell to Text format before writing, then restore the default "Normal"
# style so the visible formatting/style index is unaffected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '36.364.88'
Set-TextValue 'E2' '  -2.64%  '

# Row 3
Set-TextValue 'D3' '1.976.13'
Set-TextValue 'E3' '  -3.48%  '

# Row 4
Set-TextValue 'E4' '  +0.05%  '

# Row 5
Set-TextValue 'D5' '245.16'
Set-TextValue 'E5' '  -3.26%  '

# Row 6
Set-TextValue 'E6' '  -4.46%  '

# Row 7
Set-TextValue 'D7' '58.96'
Set-TextValue 'E7' '  -10.57%  '

# Row 8
Set-TextValue 'E8' '  +0.07%  '

# Row 9
Set-TextValue 'E9' '  -8.68%  '

# Row 10
Set-TextValue 'D10' '56.76'
Set-TextValue 'E10' '  -4.91%  '

# Row 11
Set-TextValue 'D11' '0.0854'
Set-TextValue 'E11' '  +9.04%  '

# Row 12
Set-TextValue 'E12' '  -0.47%  '

# Row 13
Set-TextValue 'D13' '22.79'
Set-TextValue 'E13' '  -3.55%  '

# Row 14
Set-TextValue 'E14' '  -7.80%  '

# Row 15
Set-TextValue 'D15' '2.265.77'
Set-TextValue 'E15' '  -3.45%  '

# Row 16
Set-TextValue 'D16' '13.82'
Set-TextValue 'E16' '  -7.11%  '

# Row 17
Set-TextValue 'D17' '5.46'
Set-TextValue 'E17' '  -5.27%  '

# Row 18
Set-TextValue 'D18' '1.976.24'
Set-TextValue 'E18' '  -3.46%  '

# Row 19
Set-TextValue 'D19' '36.328.25'
Set-TextValue 'E19' '  -2.55%  '

# Row 20
Set-TextValue 'D20' '0.0₃0884'
Set-TextValue 'E20' '  -0.46%  '

# Row 21
Set-TextValue 'D21' '70.39'
Set-TextValue 'E21' '  -4.45%  '

# Row 22
Set-TextValue 'D22' '5.27'
Set-TextValue 'E22' '  -5.28%  '

# Row 23
Set-TextValue 'D23' '233.61'
Set-TextValue 'E23' '  -2.87%  '

# Row 25
Set-TextValue 'D25' '2.51'
Set-TextValue 'E25' '  -4.87%  '

# Row 26
Set-TextValue 'E26' '  -3.05%  '

# Row 27
Set-TextValue 'D27' '9.86'
Set-TextValue 'E27' '  -3.14%  '

# Row 28
Set-TextValue 'D28' '163.25'
Set-TextValue 'E28' '  +0.68%  '

# Row 29
Set-TextValue 'E29' '  -0.63%  '

# Row 30
Set-TextValue 'D30' '19.83'
Set-TextValue 'E30' '  -1.42%  '

# Row 31
Set-TextValue 'E31' '  -2.70%  '

# Row 32
Set-TextValue 'E32' '  -1.22%  '

# Row 33
Set-TextValue 'E33' '  -6.31%  '

# Row 34
Set-TextValue 'D34' '0.0661'
Set-TextValue 'E34' '  +4.44%  '

# Row 35
Set-TextValue 'E35' '  -5.36%  '

# Row 36
Set-TextValue 'E36' '  -3.67%  '

# Row 37
Set-TextValue 'E37' '  +0.08%  '

# Row 38
Set-TextValue 'E38' '  -1.60%  '

# Row 39
Set-TextValue 'E39' '  -7.62%  '

# Row 40
Set-TextValue 'D40' '2.92'
Set-TextValue 'E40' '  -4.65%  '

# Row 41
Set-TextValue 'B41' 'TrustWalletToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D41' '1.22'
Set-TextValue 'E41' '  -4.83%  '

# Row 42
Set-TextValue 'B42' 'Cronos'
Set-TextValue 'C42' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D42' '0.0966'
Set-TextValue 'E42' '  -6.08%  '

# Row 43
Set-TextValue 'E43' '  -5.89%  '

# Row 44
Set-TextValue 'E44' '  -2.89%  '

# Row 45
Set-TextValue 'E45' '  -6.93%  '

# Row 46
Set-TextValue 'D46' '16.08'
Set-TextValue 'E46' '  -9.90%  '

# Row 47
Set-TextValue 'D47' '91.25'
Set-TextValue 'E47' '  -5.58%  '

# Row 48
Set-TextValue 'D48' '1.369.43'
Set-TextValue 'E48' '  -2.56%  '

# Row 49
Set-TextValue 'E49' '  -5.69%  '

# Row 50
Set-TextValue 'D50' '2.83'
Set-TextValue 'E50' '  -3.39%  '

# Row 51
Set-TextValue 'E51' '  -5.19%  '
